$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet holds one weekly price record per row (rows 2..134), with the
# newest record always inserted at the top of the data block (row 8 is the
# "current week" row for this particular mercado block) and every older
# record pushed down by one row. The oldest existing record (row 134) falls
# off the end of the shifted block and becomes the new last row (135).

# 1) Push the last row's data into the new row 135 first (copy, not move,
#    so row 134 keeps its values until it gets overwritten in step 2).
$ws.Range("A134:R134").Copy($ws.Range("A135:R135"))

# 2) Shift rows 9..134 down by one: row N takes what used to be in row N-1.
#    Walk bottom-up so we never clobber a source row before it has been
#    copied forward.
for ($r = 133; $r -ge 8; $r--) {
    $src = $r
    $dst = $r + 1
    $ws.Range("D$dst").Value = $ws.Range("D$src").Value2
    $ws.Range("J$dst").Value = $ws.Range("J$src").Value2
    $ws.Range("K$dst").Value = $ws.Range("K$src").Value2
    $ws.Range("L$dst").Value = $ws.Range("L$src").Value2
    $ws.Range("M$dst").Value = $ws.Range("M$src").Value2
    $ws.Range("P$dst").Value = $ws.Range("P$src").Value2
}

# 3) Row 8 becomes the new latest weekly record.
$ws.Range("D8").Value = 44496
$ws.Range("J8").Value = 60
$ws.Range("K8").Value = 6000
$ws.Range("L8").Value = 6000
$ws.Range("M8").Value = 6000
$ws.Range("P8").Value = 375
